# Update the "想去人数" (F column) counts for two rows on both the
# "展览" and "全部类型" worksheets, reflecting the latest generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 4: F4 494 -> 496
    $ws.Range("F4").Value = 496

    # Row 6: F6 6545 -> 6546
    $ws.Range("F6").Value = 6546
}
